$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values based on repulled / recalculated data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -3
$ws.Range("F9").Value = -6
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = -11
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = -11
$ws.Range("F17").Value = -2
$ws.Range("F20").Value = -7
$ws.Range("F21").Value = 0
